# Apply the Jan 4 2024 09:37:08 UTC GitHub Actions crypto-price refresh to Sheet1.
# Numeric-looking text values (e.g. plain decimals in column D) are written via
# NumberFormat '@' + Style reset so they remain Text cells (matching the source
# inlineStr cells) without leaving a stray number format behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '43.131.65'
$ws.Range('E2').Value = '  -5.12%  '
$ws.Range('D3').Value = '2.229.08'
$ws.Range('E3').Value = '  -6.18%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.582'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -8.62%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.564'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -8.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.95'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -9.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.24'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0832'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -9.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -10.11%  '
$ws.Range('E14').Value = '  -2.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.864'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -11.91%  '
$ws.Range('D16').Value = '2.562.49'
$ws.Range('E16').Value = '  -6.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.42'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.96%  '
$ws.Range('D18').Value = '2.222.23'
$ws.Range('E18').Value = '  -6.31%  '
$ws.Range('D19').Value = '43.036.60'
$ws.Range('E19').Value = '  -5.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.31%  '
$ws.Range('D21').Value = '0.0₃0969'
$ws.Range('E21').Value = '  -9.06%  '
$ws.Range('E22').Value = '  -10.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.24'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -10.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -11.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '237.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -10.35%  '
$ws.Range('E26').Value = '  -6.66%  '
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -11.11%  '
$ws.Range('E31').Value = '  -14.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '36.56'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.38'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0870'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '154.43'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.67'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.28'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.122'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.91'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.44'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.54%  '
$ws.Range('E42').Value = '  -7.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0323'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.63%  '
$ws.Range('D46').Value = '1.724.12'
$ws.Range('E46').Value = '  -8.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.14'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -13.30%  '
$ws.Range('E48').Value = '  -11.52%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.77%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.30'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -12.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '75.25'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -13.09%  '
